$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.926.01'
$ws.Range("E2").Value = '  -3.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.286.49'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '494.26'
$ws.Range("E5").Value = '  -1.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.62'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.288.04'
$ws.Range("E9").Value = '  -2.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0938'
$ws.Range("E10").Value = '  -4.10%  '
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.67'
$ws.Range("E13").Value = '  -3.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.691.71'
$ws.Range("E14").Value = '  -2.10%  '
$ws.Range("E15").Value = '  +0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '53.935.05'
$ws.Range("E16").Value = '  -3.07%  '
$ws.Range("E17").Value = '  -2.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.280.15'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.91'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.04'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '298.91'
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").Value = '  +1.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -2.61%  '
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("E26").Value = '  +1.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.352.84'
$ws.Range("E27").Value = '  -3.84%  '
$ws.Range("E28").Value = '  +1.66%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '163.41'
$ws.Range("E30").Value = '  -5.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.60'
$ws.Range("E31").Value = '  -1.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0₃0681'
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.84'
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.871'
$ws.Range("E39").Value = '  +5.50%  '
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.37'
$ws.Range("E41").Value = '  -1.94%  '
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("E45").Value = '  +4.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.38'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("E48").Value = '  -1.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '238.35'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0480'
$ws.Range("E50").Value = '  +0.59%  '
$ws.Range("E51").Value = '  -0.63%  '
